# Generate Report for Archive
# - Update the localization status from "Ready for handoff" to "In Translation"
#   on the Overview sheet (columns E/F, row 2) and on each language sheet
#   (zh-cn, de-de) in the Status column (C, row 2).
# - Re-fit the now-narrower Status columns so their stored width matches the
#   shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# Overview sheet: zh-cn / de-de status columns are E and F on row 2.
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Per-language detail sheets: Status column is C on row 2.
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# Narrow the Status columns to fit the new, shorter text.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
